# Convert the plain Heading1/bold-run title block into a pandoc-style
# title block: a "Title"-styled paragraph for the headline and an
# "Authors"-styled paragraph for the byline, each split into one run
# per word/separator (matching the pandoc docx writer's output shape).

$d = $word.ActiveDocument
$wns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# --- Paragraph 1: "On Pilgrimage - July/August 1975" (Heading1 -> Title) ---
$titleXml = @"
<w:p $wns><w:pPr><w:pStyle w:val="Title"/></w:pPr><w:r><w:t xml:space="preserve">On</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Pilgrimage</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">-</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">July</w:t></w:r><w:r><w:t xml:space="preserve">/</w:t></w:r><w:r><w:t xml:space="preserve">August</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">1975</w:t></w:r></w:p>
"@

$d.Paragraphs(1).Range.InsertXML($titleXml) | Out-Null

# --- Paragraph 2: "By Dorothy Day" (bold run -> Authors style, "Dorothy Day") ---
$authorsXml = @"
<w:p $wns><w:pPr><w:pStyle w:val="Authors"/></w:pPr><w:r><w:t xml:space="preserve">Dorothy</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Day</w:t></w:r></w:p>
"@

$d.Paragraphs(2).Range.InsertXML($authorsXml) | Out-Null
